$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 and C3 values
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 9120044824

# Set column C width (closest achievable value to the target 17.7109375 chars
# given this runtime's internal pixel-rounding grid)
$ws.Columns.Item(3).ColumnWidth = 16.83

# Change the active selection to C2
$ws.Range("C2").Select()
